$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the existing "sum" header (G1) onto the new H1 header
# cell so it reuses the same bold/border/center-top style used by the other
# header cells, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New "Save" data column value for row 2
$ws.Range("H2").Value = 0
